$d = $word.ActiveDocument

# --- Change 1: "Introduction to Advance Smart GPS System" paragraph ---
# "Introduction to " -> "Introduction and " + "History" (rFonts cstheme=minorHAnsi) + " of "
$r1 = $d.Content
$r1.Find.Execute("Introduction to Advance Smart GPS System")
$p1 = $r1.Paragraphs(1).Range
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="007C036D" w:rsidRPr="00040D55" w:rsidRDefault="007C036D" w:rsidP="007C036D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:jc w:val="both"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Introduction and </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>History</w:t></w:r><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:r w:rsidRPr="00C40AA7"><w:t>Advance Smart GPS System</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p1.InsertXML($xml1)

# --- Change 2: "Introduction to Electromagnetic Brake system" paragraph ---
# " to " -> " and " + "History" (rFonts cstheme=minorHAnsi) + " of " (all rFonts cstheme=minorHAnsi)
$r2 = $d.Content
$r2.Find.Execute("Introduction to Electromagnetic Brake system")
$p2 = $r2.Paragraphs(1).Range
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="007C036D" w:rsidRPr="00010BB3" w:rsidRDefault="007C036D" w:rsidP="007C036D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:jc w:val="both"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r w:rsidRPr="003E6DDD"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>Introduction</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>History</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> of </w:t></w:r><w:r w:rsidRPr="00A24A08"><w:t>Electromagnetic Brake system</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.InsertXML($xml2)

# --- Change 3: "History and Literature Review of " -> "Literature Review of " (Advance Smart GPS System heading) ---
$d.Content.Find.Execute("History and Literature Review of ", $false, $false, $false, $false, $false, $true, 1, $false, "Literature Review of ", 2)

# --- Change 4: "History and Literature Review of Electromagnetic Brake System" -> "Literature Review of Electromagnetic Brake System" ---
$d.Content.Find.Execute("History and Literature Review of Electromagnetic Brake System", $false, $false, $false, $false, $false, $true, 1, $false, "Literature Review of Electromagnetic Brake System", 2)

Write-Output "edit complete"
